$d = $word.ActiveDocument

# The paragraph currently reads (single run):
#   "   - Menu posting interface: Cooks can upload their daily menu (with
#    descriptions, price, and availability)."
# We need to keep the leading three spaces unformatted, and apply a green
# highlight to the "- Menu posting interface..." portion only. Doing this
# via Find on just that substring (excluding the leading spaces) causes
# Word to split the run automatically at the point where the new
# character formatting starts/ends.
$needle = "- Menu posting interface: Cooks can upload their daily menu (with descriptions, price, and availability)."

$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if (-not $found) {
    throw "Could not find target text: $needle"
}

# wdBrightGreen (4) is the WdColorIndex value that round-trips to the
# OOXML <w:highlight w:val="green"/> element. Setting it through the
# Range's Font object (rather than directly on the Range) correctly
# scopes the formatting change to just this Range/Find hit.
$rng.Font.HighlightColorIndex = 4
